# Chapter 5 additions to the glossary word list:
#   - fundamental theorem of integral calculus / 적분학의 기본정리
#   - maximum principle / 최대원리
# Both are tagged with the "용어사전" (glossary) note in column D,
# appended as new rows right after the existing last row (101).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 102
$ws.Cells.Item(102, 1).Value = "fundamental theorem of integral calculus"
$ws.Cells.Item(102, 2).Value = "적분학의 기본정리"
$ws.Cells.Item(102, 4).Value = "용어사전"

# New row 103
$ws.Cells.Item(103, 1).Value = "maximum principle"
$ws.Cells.Item(103, 2).Value = "최대원리"
$ws.Cells.Item(103, 4).Value = "용어사전"

# Move the view/selection forward to the next empty row, mirroring the
# author's cursor position after typing the new entries.
$win = $excel.ActiveWindow
$win.ScrollRow = 92
$win.ScrollColumn = 1
$ws.Range("A104").Select() | Out-Null
